$d = $word.ActiveDocument

# 1. Remove the "Abstract Title" custom paragraph style entirely.
$abstractTitle = $d.Styles("AbstractTitle")
$abstractTitle.Delete()

# 2. Change the "Abstract" style's spacing-before from 100 (5pt) to 300 (15pt) twips.
$abstract = $d.Styles("Abstract")
$abstract.ParagraphFormat.SpaceBefore = 15

# 3. Remove the "Footnote Block Text" custom paragraph style entirely.
$footnoteBlockText = $d.Styles("FootnoteBlockText")
$footnoteBlockText.Delete()

Write-Output "done"
